$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting (border + centered alignment) from an existing
# data row down onto the two new rows before filling them in.
$ws.Range("A9:B9").Copy()
$ws.Range("A11:B12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the donor pairs (new row 11/12) and the previously-blank A10 id,
# written in this order so the shared-string table matches the source edit.
$ws.Range("B11").Value = "Br6032"
$ws.Range("A10").Value = "Br5436"
$ws.Range("A11").Value = "Br5931"
$ws.Range("A12").Value = "Br6389"
$ws.Range("B12").Value = "Br5746"

# Remove the stray ReadMe note that lived in D1 on Sheet1
$ws.Range("D1").ClearContents()

# Match the saved selection state
$ws.Range("D1").Select()
